$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Hunk 1: Main scenario step "7" -> "5" ("Sistema: Controlla che i dati inseriti ...")
$c = $t.Rows.Item(20).Cells.Item(1).Range
$c.Find.Execute("7", $false, $false, $false, $false, $false, $true, 1, $false, "5", 1)

# Hunk 2: Main scenario step "8" -> "6" ("Sistema: Convalida i dati di accesso ...")
$c = $t.Rows.Item(21).Cells.Item(1).Range
$c.Find.Execute("8", $false, $false, $false, $false, $false, $true, 1, $false, "6", 1)

# Hunk 3: "I" + "I" + " Scenario/Flusso di eventi Alternativo:  " runs merge into
# a single "II Scenario/Flusso di eventi Alternativo:  " run (same visible text,
# re-saving the range coalesces the three runs that share identical formatting).
$c = $t.Rows.Item(27).Cells.Item(1).Range
$c.Find.Execute("II Scenario/Flusso di eventi Alternativo:", $true, $false, $false, $false, $false, $true, 1, $false, "II Scenario/Flusso di eventi Alternativo:", 1)

# Hunk 4: II Scenario step "8.a1" -> "5.a1" (only the leading "8" run changes)
$c = $t.Rows.Item(28).Cells.Item(1).Range
$c.Find.Execute("8", $false, $false, $false, $false, $false, $true, 1, $false, "5", 1)

# Hunk 5: II Scenario step "8.a2" -> "5.a2" (only the leading "8" run changes)
$c = $t.Rows.Item(29).Cells.Item(1).Range
$c.Find.Execute("8", $false, $false, $false, $false, $false, $true, 1, $false, "5", 1)

# Hunk 6: Error scenario step "7.a1" -> "5.a1" (only the leading "7" run changes)
$c = $t.Rows.Item(32).Cells.Item(1).Range
$c.Find.Execute("7", $false, $false, $false, $false, $false, $true, 1, $false, "5", 1)

# Hunk 7: Error scenario step "7.a2" -> "5.a2" (only the leading "7" run changes)
$c = $t.Rows.Item(33).Cells.Item(1).Range
$c.Find.Execute("7", $false, $false, $false, $false, $false, $true, 1, $false, "5", 1)
